$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.085.02'
$ws.Range('E2').Value = '  -1.77%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.891.24'
$ws.Range('E3').Value = '  -1.15%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.27'
$ws.Range('E5').Value = '  -0.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.07%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5028'
$ws.Range('E7').Value = '  -1.81%  '

$ws.Range('E8').Value = '  -1.73%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09175'
$ws.Range('E9').Value = '  -6.09%  '

$ws.Range('E10').Value = '  -3.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.88'
$ws.Range('E11').Value = '  -0.74%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.378'

$ws.Range('E13').Value = '  -2.18%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.895.32'
$ws.Range('E14').Value = '  -1.03%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.276'
$ws.Range('E15').Value = '  -4.18%  '

$ws.Range('E16').Value = '  +0.19%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.37'
$ws.Range('E17').Value = '  -1.74%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001104'
$ws.Range('E18').Value = '  -3.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06671'
$ws.Range('E19').Value = '  +0.08%  '

$ws.Range('E20').Value = '  -2.21%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.197'
$ws.Range('E22').Value = '  -2.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.132.24'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.37'
$ws.Range('E24').Value = '  -0.94%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.320'
$ws.Range('E25').Value = '  +1.40%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.117.76'
$ws.Range('E26').Value = '  -0.75%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.533'
$ws.Range('E27').Value = '  -7.25%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '158.03'
$ws.Range('E28').Value = '  -1.01%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.81'
$ws.Range('E29').Value = '  -2.51%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.86'
$ws.Range('E30').Value = '  -1.56%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.073'
$ws.Range('E31').Value = '  -2.87%  '

$ws.Range('E32').Value = '  -2.18%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.596'
$ws.Range('E33').Value = '  -2.73%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.620'
$ws.Range('E34').Value = '  -0.71%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.468'
$ws.Range('E35').Value = '  -4.41%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06574'
$ws.Range('E36').Value = '  -3.59%  '

$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.337'
$ws.Range('E37').Value = '  +12.08%  '

$ws.Range('E38').Value = '  -2.02%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2200'
$ws.Range('E39').Value = '  -1.63%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.212'
$ws.Range('E40').Value = '  -4.76%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6432'
$ws.Range('E41').Value = '  -0.45%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.43'
$ws.Range('E42').Value = '  -4.46%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.954'
$ws.Range('E43').Value = '  -3.08%  '

$ws.Range('E44').Value = '  +0.14%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.26'
$ws.Range('E45').Value = '  -2.98%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6040'
$ws.Range('E46').Value = '  -1.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.303'
$ws.Range('E47').Value = '  +1.62%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.684'

$ws.Range('E49').Value = '  -2.56%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '122.23'
$ws.Range('E50').Value = '  -2.56%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.193'
$ws.Range('E51').Value = '  -1.83%  '
